$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Wed Dec 20 12:57:28 EST 2023"
$ws.Range("B3").Value = "Wed Dec 20 12:57:40 EST 2023"
$ws.Range("B5").Value = "Wed Dec 20 12:57:52 EST 2023"
